$d = $word.ActiveDocument

# Helper constants
$wdFindStop = 0
$wdReplaceNone = 0
$wdReplaceOne = 1
$wdReplaceAll = 2

function Underline-FirstMatch($paraIndex, $searchText) {
    $r = $d.Paragraphs.Item($paraIndex).Range
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindStop, $false, "", $wdReplaceNone)
    if ($found) {
        $r.Font.Underline = 1
    }
    return $found
}

function ReplaceText-InParagraph($paraIndex, $searchText, $replaceText) {
    $r = $d.Paragraphs.Item($paraIndex).Range
    $found = $r.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, $wdFindStop, $false, $replaceText, $wdReplaceOne)
    return $found
}

# 1) The "Costos" chart drawing run: lang es-AR -> noProof + es-ES_tradnl
$costosPara = $d.Paragraphs.Item(109 - 109)  # placeholder, replaced below
$found = $d.Content.Find.Execute("Costos", $true, $false, $false, $false, $false, $true, $wdFindStop, $false, "", $wdReplaceNone)
Write-Output ("found Costos=" + $found)

# 2) CERRADOS section - "Mitigación: Se acordó con Alejandro" (para 91)
Underline-FirstMatch 91 "Mitigación:" | Out-Null

# 3) "Fecha de cierre: Reunión Formal..." (para 92)
Underline-FirstMatch 92 "Fecha de cierre:" | Out-Null

# 4) "Exposición: 0,20 (Disminución ... riguroso)" (para 100)
Underline-FirstMatch 100 "Exposición:" | Out-Null

# 5) "Mitigación: Comenzar lo antes posible..." (para 101)
Underline-FirstMatch 101 "Mitigación:" | Out-Null

# 6) "Exposición" + ": 0,2" (para 105) -> underline "Exposición:"
Underline-FirstMatch 105 "Exposición:" | Out-Null

# 7) "Mitigación" + ": Hablar..." (para 106) -> underline "Mitigación:"
Underline-FirstMatch 106 "Mitigación:" | Out-Null

# 8) "Recepción de datos incorrectos." -> "Aspectos del negocio poco claros por parte del equipo" (para 108)
ReplaceText-InParagraph 108 "Recepción de datos incorrectos." "Aspectos del negocio poco claros por parte del equipo" | Out-Null

# 9) "Exposición: 0,20 " -> underline "Exposición:" and change value to 0,15 (para 109)
ReplaceText-InParagraph 109 "0,20" "0,15" | Out-Null
Underline-FirstMatch 109 "Exposición:" | Out-Null

# 10) "Mitigación: Establecer una capa intermedia..." -> "Mitigación: Hacer consultas al cliente sobre el negocio de forma constante." (para 110)
ReplaceText-InParagraph 110 "Establecer una capa intermedia de validación de datos y con reporte de warnings y errores encontrados." "Hacer consultas al cliente sobre el negocio de forma constante." | Out-Null
Underline-FirstMatch 110 "Mitigación:" | Out-Null

Write-Output "done"
